$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CIES-11815", "Prioridade 3", "CÂMERAS DO P0214 OBSTRUIDAS", "Service Desk Perkons", "P0214-PK-VLS-JAG", "FORA DO PRAZO", "NO PRAZO"),
    @("CIES-11821", "Prioridade 1", "todas as câmeras do p0173 estão off-line", "CRISTIANA NEVES BATISTA MOREIRA", "P0173-PK-VLS-IBA", "NO PRAZO", "NO PRAZO"),
    @("CIES-11823", "Prioridade 1", "todas as câmeras do p0271 estão off-line", "CRISTIANA NEVES BATISTA MOREIRA", "P0271-PCL-PK-ARA", "NO PRAZO", "NO PRAZO"),
    @("CIES-11825", "Prioridade 1", "todas as câmeras do p0271 estão off-line", "CRISTIANA NEVES BATISTA MOREIRA", "P0259-PCL-PK-COL", "NO PRAZO", "NO PRAZO"),
    @("CIES-11827", "Prioridade 1", "Todas as câmeras do p0234 estão off-line", "CRISTIANA NEVES BATISTA MOREIRA", "P0234-PK-VLS-PED", "NO PRAZO", "NO PRAZO"),
    @("CIES-11828", "Prioridade 1", "Todas as câmeras do p0243 estão off-line", "CRISTIANA NEVES BATISTA MOREIRA", "P0243-PK-VLS-NOV", "NO PRAZO", "NO PRAZO"),
    @("CIES-11829", "Prioridade 1", "Todas as câmeras do p0229 estão off-line", "CRISTIANA NEVES BATISTA MOREIRA", "P0229-PK-VLS-VPA", "NO PRAZO", "NO PRAZO"),
    @("CIES-11834", "Prioridade 1", "todas as câmeras do p0249 estão off-line  ", "CRISTIANA NEVES BATISTA MOREIRA", "P0249-PK-VLS-STE", "NO PRAZO", "NO PRAZO"),
    @("CIES-11835", "Prioridade 1", "câmeras do p0293 estão off-line ", "CRISTIANA NEVES BATISTA MOREIRA", "P0293-PK-VLS-CAC", "NO PRAZO", "NO PRAZO"),
    @("CIES-12146", "Prioridade 1", "Corretiva P0174", "Service Desk Perkons", "P0174-PK-VLS-LIN", "NO PRAZO", "NO PRAZO"),
    @("CIES-12149", "Prioridade 2", " câmeras do p0272 estão off-line  ", "CRISTIANA NEVES BATISTA MOREIRA", "P0272-PCL-PK-ALF", "NO PRAZO", "NO PRAZO"),
    @("CIES-12151", "Prioridade 1", "P0281 - Site Offline", "Service Desk Perkons", "P0281-PCL-PK-PAN", "NO PRAZO", "NO PRAZO"),
    @("CIES-12152", "Prioridade 1", "P0243 - Site sem acesso externo", "Service Desk Perkons", "P0243-PK-VLS-NOV", "NO PRAZO", "NO PRAZO"),
    @("CIES-12153", "Prioridade 1", "CORRETIVA SITE P0281", "Service Desk Perkons", "P0281-PCL PK-PAN", "NO PRAZO", "NO PRAZO"),
    @("CIES-12155", "Prioridade 1", "CORRETIVA SITE P0243", "Service Desk Perkons", "P0243-PK-VLS-NOV", "NO PRAZO", "NO PRAZO"),
    @("CIES-12157", "Prioridade 1", "todas as câmeras do p0226 estão off-line", "CRISTIANA NEVES BATISTA MOREIRA", "P0226-PK-VLS-JER", "NO PRAZO", "NO PRAZO"),
    @("CIES-12160", "Prioridade 1", "CORRETIVA SITE P0184", "Service Desk Perkons", "P0184-PK-VLS-GUA", "NO PRAZO", "NO PRAZO"),
    @("CIES-12162", "Prioridade 1", "CORRETIVA SITE P0215", "Service Desk Perkons", "P0215-PK-VLS-NOV", "NO PRAZO", "NO PRAZO"),
    @("CIES-12163", "Prioridade 1", "CORRETIVA SITE P0249", "Service Desk Perkons", "P0249-PK-VLS-STE", "NO PRAZO", "NO PRAZO"),
    @("CIES-12171", "Prioridade 1", "CORRETIVA SITE P0259", "Service Desk Perkons", "P0259-PCL-PK-COL", "NO PRAZO", "NO PRAZO"),
    @("CIES-12172", "Prioridade 1", "CORRETIVA SITE P0308", "Service Desk Perkons", "P0308-PCL-PK-DOM", "NO PRAZO", "NO PRAZO"),
    @("CIES-12174", "Prioridade 1", "P0205 - CORRETIVA SITE", "Service Desk Perkons", "P0205-PK-VLS-BAI", "NO PRAZO", "NO PRAZO"),
    @("CIES-12179", "Prioridade 1", "todas as câmeras do p0306 estão off-line   ", "CRISTIANA NEVES BATISTA MOREIRA", "P0306-PCL-PK-MIM", "NO PRAZO", "NO PRAZO"),
    @("CIES-12180", "Prioridade 1", "todas as câmeras do p0303 estão off-line   ", "CRISTIANA NEVES BATISTA MOREIRA", "P0303-PK-VLS-IRU", "NO PRAZO", "NO PRAZO"),
    @("CIES-12181", "Prioridade 1", "todas as câmeras do p0247 estão off-line   ", "CRISTIANA NEVES BATISTA MOREIRA", "P0247-PK-VLS-ADN", "NO PRAZO", "NO PRAZO"),
    @("CIES-12183", "Prioridade 1", "todas as câmeras do p0281 estão off-line ", "CRISTIANA NEVES BATISTA MOREIRA", "P0281-PCL PK-PAN", "NO PRAZO", "NO PRAZO"),
    @("CIES-12184", "Prioridade 1", "todas as câmeras do p0303 estão off-line ", "CRISTIANA NEVES BATISTA MOREIRA", "P0203-PK-VLS-ECO", "NO PRAZO", "NO PRAZO"),
    @("CIES-12218", "Prioridade 1", "CORRETIVA SITE P0249", "Service Desk Perkons", "P0249-PK-VLS-STE", "NO PRAZO", "NO PRAZO"),
    @("CIES-12219", "Prioridade 1", "CORRETIVA SITE P0259", "Service Desk Perkons", "P0259-PCL-PK-COL", "NO PRAZO", "NO PRAZO"),
    @("CIES-12229", "Prioridade 1", "CORRETIVA SITE P0205", "Service Desk Perkons", "P0205-PK-VLS-BAI", "NO PRAZO", "NO PRAZO"),
    @("CIES-12231", "Prioridade 1", "P0215 - Todas câmeras offline. ", "Service Desk Perkons", "P0215-PK-VLS-NOV", "NO PRAZO", "NO PRAZO"),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Host "Done applying data"
